# Facilitators guidelines - Surface Tension: English -> Swahili (Kenya) translation
$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Longer / more specific phrases first so they are not clobbered by the
# shorter phrases that are substrings of them (e.g. "Video Introduction"
# is contained in "General VMC Video Introduction").
Replace-Text "General VMC Video Introduction" "Utangulizi Mkuu wa Video ya VMC"
Replace-Text "Video Introduction" "Utangulizi wa Video"

Replace-Text "Video Title" "Kichwa cha Video"
Replace-Text "Topic" "Mada"
Replace-Text "Aim(s)" "Malengo"
Replace-Text "Length" "Urefu"
Replace-Text "Camp Location" "Mahali pa Kambi"
Replace-Text "Facilitators" "Wawezeshaji"
Replace-Text "N. of students" "N. ya wanafunzi"
Replace-Text "Date" "Tarehe"
Replace-Text "Resources" "Rasilimali"
Replace-Text "needed" "inahitajika"
Replace-Text "Preparations" "Maandalizi"
Replace-Text "Video time" "Muda wa video"
Replace-Text "What facilitator does" "Mwezeshaji anafanya nini"
Replace-Text "What learners do" "Wanachofanya wanafunzi"
Replace-Text "Introduction of the first experiment" "Utangulizi wa jaribio la kwanza"
# Occurs twice in the document; ReplaceAll (wdReplaceAll = 2) handles both.
Replace-Text "Assist the process, provoke thoughts" "Kusaidia mchakato, kuchochea mawazo"

# Document default language: Swahili (Tanzania) -> Swahili (Kenya)
$d.Styles.Item("Normal").LanguageID = "sw-KE"
